$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before column B (old B,C,D shift right to E,F,G)
$ws.Range("B1:D1").EntireColumn.Insert()

# Header row for new columns
$ws.Range("B1").Value = "analysis_results_20250501"
$ws.Range("C1").Value = "analysis_results_20250430"
$ws.Range("D1").Value = "analysis_results_20250429"

# Copy header style from the (now shifted) original header cell E1 onto the new headers
$ws.Range("E1").Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Values for each data row: B=20250501 result, C=20250430 result, D=20250429 result
$results = @{
  2 = "error"
  3 = "success"
  4 = "success"
  5 = "success"
  6 = "error"
  7 = "success"
  8 = "success"
  9 = "error"
  10 = "success"
  11 = "error"
  12 = "success"
  13 = "success"
}

foreach ($row in 2..13) {
  $val = $results[$row]
  $ws.Range("B$row").Value = $val
  $ws.Range("C$row").Value = $val
  $ws.Range("D$row").Value = $val
}
